$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-formatted cells (values that would otherwise be auto-parsed as numbers)
# to retain their exact original text representation, matching the source data feed.
$textCells = @("D5","D8","D10","D11","D13","D15","D16","D20","D22","D23","D24","D26","D30","D32","D33","D34","D35","D40","D42","D46","D47","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values scraped from the latest cryptos feed
$ws.Range("D2").Value = "36.561.86"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "1.995.67"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "233.76"
$ws.Range("E5").Value = "  -9.39%  "
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "55.08"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").Value = "57.89"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "14.24"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.288.92"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "20.42"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "0.759"
$ws.Range("E16").Value = "  -5.11%  "
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "2.003.16"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "36.462.40"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").Value = "67.77"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "0.0₃0805"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").Value = "5.30"
$ws.Range("E22").Value = "  +3.39%  "
$ws.Range("D23").Value = "222.14"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -8.47%  "
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "18.81"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "4.37"
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("D34").Value = "0.0603"
$ws.Range("E34").Value = "  -5.98%  "
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").Value = "  -6.27%  "
$ws.Range("E36").Value = "  -1.05%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  -3.34%  "
$ws.Range("D40").Value = "5.64"
$ws.Range("E40").Value = "  +7.62%  "
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "0.0945"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("D43").Value = "1.454.97"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("E44").Value = "  -4.33%  "
$ws.Range("E45").Value = "  -8.87%  "
$ws.Range("D46").Value = "89.12"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "15.22"
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("E48").Value = "  -2.69%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("D51").Value = "3.73"
$ws.Range("E51").Value = "  +8.14%  "

Write-Output "Applied cryptos list update"
